$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.001362816872990891
$ws.Range("E2").Value = 0.4336984297763564
$ws.Range("F2").Value = 0.716057570907509
$ws.Range("G2").Value = 0.002324725458065281
$ws.Range("O2").Value = 2.355210718842727

$ws.Range("D3").Value = 0.001187382742548593
$ws.Range("E3").Value = 0.377859479011633
$ws.Range("F3").Value = 0.6567163053089047
$ws.Range("G3").Value = 0.002329461130132879
$ws.Range("O3").Value = 2.156718067769077

$ws.Range("D4").Value = 0.001080080273940354
$ws.Range("E4").Value = 0.3437100912600926
$ws.Range("F4").Value = 0.6207120828670298
$ws.Range("G4").Value = 0.00233251953125169
$ws.Range("O4").Value = 2.036265997455303

$ws.Range("D5").Value = 0.001036432188708858
$ws.Range("E5").Value = 0.3298250761745436
$ws.Range("F5").Value = 0.6061468742205847
$ws.Range("G5").Value = 0.002333803886739488
$ws.Range("O5").Value = 1.987532592803802

$ws.Range("D6").Value = 0.00102918845908917
$ws.Range("E6").Value = 0.3275212721867717
$ws.Range("F6").Value = 0.6037347334826109
$ws.Range("G6").Value = 0.002334019454191858
$ws.Range("O6").Value = 1.979461519004019

$ws.Range("D7").Value = 0.001079491336421157
$ws.Range("E7").Value = 0.3435227109891485
$ws.Range("F7").Value = 0.6205152211821741
$ws.Range("G7").Value = 0.002332536698358002
$ws.Range("O7").Value = 2.035607345351195

$ws.Range("D8").Value = 0.001302225394198331
$ws.Range("E8").Value = 0.4144150822295103
$ws.Range("F8").Value = 0.6955060775011219
$ws.Range("G8").Value = 0.002326327130478486
$ws.Range("O8").Value = 2.286471316272923

$ws.Range("D9").Value = 0.001743495442950049
$ws.Range("E9").Value = 0.5546583365364199
$ws.Range("F9").Value = 0.8460728912957336
$ws.Range("G9").Value = 0.00231533923825636
$ws.Range("O9").Value = 2.79001036887513

$ws.Range("D10").Value = 0.002072145766176092
$ws.Range("E10").Value = 0.6586520499176061
$ws.Range("F10").Value = 0.9589651656520743
$ws.Range("G10").Value = 0.002307982165639402
$ws.Range("O10").Value = 3.167490004210549

$ws.Range("D11").Value = 0.002223006804682015
$ws.Range("E11").Value = 0.7062161564137313
$ws.Range("F11").Value = 1.010845266279489
$ws.Range("G11").Value = 0.002304788697645446
$ws.Range("O11").Value = 3.340953389518347

$ws.Range("D12").Value = 0.002280360919042934
$ws.Range("E12").Value = 0.7242680764979212
$ws.Range("F12").Value = 1.030568565905213
$ws.Range("G12").Value = 0.00230360130796144
$ws.Range("O12").Value = 3.406898250946426

$ws.Range("D13").Value = 0.002267998111644687
$ws.Range("E13").Value = 0.7203784183419941
$ws.Range("F13").Value = 1.026317325067708
$ws.Range("G13").Value = 0.002303856061437034
$ws.Range("O13").Value = 3.392684256368284

$ws.Range("D14").Value = 0.002227720644000897
$ws.Range("E14").Value = 0.7077004683360428
$ws.Range("F14").Value = 1.012466350817789
$ws.Range("G14").Value = 0.002304690572135448
$ws.Range("O14").Value = 3.346373501059304

$ws.Range("D15").Value = 0.002203079946879072
$ws.Range("E15").Value = 0.6999402271471666
$ws.Range("F15").Value = 1.003992371776604
$ws.Range("G15").Value = 0.002305204582881486
$ws.Range("O15").Value = 3.31804063982969

$ws.Range("D16").Value = 0.002062316590965452
$ws.Range("E16").Value = 0.6555490890936539
$ws.Range("F16").Value = 0.9555854069585621
$ws.Range("G16").Value = 0.002308193939339283
$ws.Range("O16").Value = 3.15618949106215

$ws.Range("D17").Value = 0.00197633221307747
$ws.Range("E17").Value = 0.6283848123191689
$ws.Range("F17").Value = 0.9260251098173171
$ws.Range("G17").Value = 0.00231006697917523
$ws.Range("O17").Value = 3.057351167739625

$ws.Range("D18").Value = 0.001927002019428414
$ws.Range("E18").Value = 0.6127846016899667
$ws.Range("F18").Value = 0.9090720897889355
$ws.Range("G18").Value = 0.002311158739864443
$ws.Range("O18").Value = 3.000665896692567

$ws.Range("D19").Value = 0.001910320321343306
$ws.Range("E19").Value = 0.6075066487270959
$ws.Range("F19").Value = 0.9033404933890523
$ws.Range("G19").Value = 0.002311530874994706
$ws.Range("O19").Value = 2.981501180670932

$ws.Range("D20").Value = 0.001985472159223178
$ws.Range("E20").Value = 0.6312739907923657
$ws.Range("F20").Value = 0.9291667394193155
$ws.Range("G20").Value = 0.002309866097719104
$ws.Range("O20").Value = 3.067855660969713

$ws.Range("D21").Value = 0.002239544717282271
$ws.Range("E21").Value = 0.7114231643078313
$ws.Range("F21").Value = 1.016532600446681
$ws.Range("G21").Value = 0.002304444862416544
$ws.Range("O21").Value = 3.359969033151742

$ws.Range("D22").Value = 0.002406930070236513
$ws.Range("E22").Value = 0.7640425772890183
$ws.Range("F22").Value = 1.074083505283085
$ws.Range("G22").Value = 0.002301029401397889
$ws.Range("O22").Value = 3.552389396738079

$ws.Range("D23").Value = 0.002317460720473008
$ws.Range("E23").Value = 0.7359357237643849
$ws.Range("F23").Value = 1.043325465387227
$ws.Range("G23").Value = 0.002302840664476058
$ws.Range("O23").Value = 3.449550766765867

$ws.Range("D24").Value = 0.001981339669695714
$ws.Range("E24").Value = 0.6299677408470501
$ws.Range("F24").Value = 0.9277462792848041
$ws.Range("G24").Value = 0.002309956869656864
$ws.Range("O24").Value = 3.06310614965048

$ws.Range("D25").Value = 0.001623465342587949
$ws.Range("E25").Value = 0.5165662124021964
$ws.Range("F25").Value = 0.804951376280016
$ws.Range("G25").Value = 0.002318185404698143
$ws.Range("O25").Value = 2.652501122830699

